# Updates country rows (Benin/Etiopia/Montserrat reordering) and refreshes
# the covid numeric stats + 'last updated' timestamp for paises.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1) - data refresh time moved from 12:04 to 13:09
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 8 de Mayo de 2020 a las 13:09"

$ws.Cells.Item(13, 2).Value = 104691
$ws.Cells.Item(13, 3).Value = 1556
$ws.Cells.Item(13, 4).Value = 83837
$ws.Cells.Item(13, 5).Value = 14313
$ws.Cells.Item(13, 6).Value = 2711
$ws.Cells.Item(13, 7).Value = 55
$ws.Cells.Item(13, 8).Value = 6541

$ws.Cells.Item(52, 2).Value = 6914
$ws.Cells.Item(52, 3).Value = 18
$ws.Cells.Item(52, 4).Value = 6079
$ws.Cells.Item(52, 5).Value = 738
$ws.Cells.Item(52, 6).Value = 21
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 97

$ws.Cells.Item(56, 2).Value = 5661
$ws.Cells.Item(56, 3).Value = 113
$ws.Cells.Item(56, 4).Value = 2302
$ws.Cells.Item(56, 5).Value = 3174
$ws.Cells.Item(56, 6).Value = 1
$ws.Cells.Item(56, 7).Value = 2
$ws.Cells.Item(56, 8).Value = 185

$ws.Cells.Item(77, 2).Value = 2070
$ws.Cells.Item(77, 3).Value = 43
$ws.Cells.Item(77, 4).Value = 960
$ws.Cells.Item(77, 5).Value = 1012
$ws.Cells.Item(77, 6).Value = 4
$ws.Cells.Item(77, 7).Value = 8
$ws.Cells.Item(77, 8).Value = 98

$ws.Cells.Item(104, 2).Value = 796
$ws.Cells.Item(104, 3).Value = 12
$ws.Cells.Item(104, 4).Value = 223
$ws.Cells.Item(104, 5).Value = 547
$ws.Cells.Item(104, 6).Value = 42
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = 26

$ws.Cells.Item(133, 2).Value = 288
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(133, 4).Value = 241
$ws.Cells.Item(133, 5).Value = 47
$ws.Cells.Item(133, 6).Value = 8
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 0

$ws.Cells.Item(137, 1).Value = "Benin"
$ws.Cells.Item(137, 2).Value = 242
$ws.Cells.Item(137, 3).Value = 102
$ws.Cells.Item(137, 4).Value = 62
$ws.Cells.Item(137, 5).Value = 178
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 2

$ws.Cells.Item(138, 1).Value = "Sierra Leona"
$ws.Cells.Item(138, 2).Value = 231
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 54
$ws.Cells.Item(138, 5).Value = 161
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 16

$ws.Cells.Item(139, 1).Value = "Cabo Verde"
$ws.Cells.Item(139, 2).Value = 218
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(139, 4).Value = 38
$ws.Cells.Item(139, 5).Value = 178
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 2

$ws.Cells.Item(140, 1).Value = "Etiopia"
$ws.Cells.Item(140, 2).Value = 194
$ws.Cells.Item(140, 3).Value = 3
$ws.Cells.Item(140, 4).Value = 95
$ws.Cells.Item(140, 5).Value = 95
$ws.Cells.Item(140, 6).Value = 1
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 4

$ws.Cells.Item(141, 1).Value = "Madagascar"
$ws.Cells.Item(141, 2).Value = 193
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 4).Value = 101
$ws.Cells.Item(141, 5).Value = 92
$ws.Cells.Item(141, 6).Value = 1
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 0

$ws.Cells.Item(142, 1).Value = "Liberia"
$ws.Cells.Item(142, 2).Value = 189
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(142, 4).Value = 79
$ws.Cells.Item(142, 5).Value = 90
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 20

$ws.Cells.Item(143, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(143, 2).Value = 187
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 4
$ws.Cells.Item(143, 5).Value = 179
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 4

$ws.Cells.Item(144, 1).Value = "Islas Feroe"
$ws.Cells.Item(144, 2).Value = 187
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 186
$ws.Cells.Item(144, 5).Value = 1
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 0

$ws.Cells.Item(145, 1).Value = "Martinica"
$ws.Cells.Item(145, 2).Value = 183
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 83
$ws.Cells.Item(145, 5).Value = 86
$ws.Cells.Item(145, 6).Value = 3
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 14

$ws.Cells.Item(146, 1).Value = "Birmania"
$ws.Cells.Item(146, 2).Value = 176
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 62
$ws.Cells.Item(146, 5).Value = 108
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 6

$ws.Cells.Item(147, 1).Value = "Suazilandia"
$ws.Cells.Item(147, 2).Value = 153
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 12
$ws.Cells.Item(147, 5).Value = 139
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 2

$ws.Cells.Item(148, 1).Value = "Zambia"
$ws.Cells.Item(148, 2).Value = 153
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 103
$ws.Cells.Item(148, 5).Value = 46
$ws.Cells.Item(148, 6).Value = 1
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 4

$ws.Cells.Item(149, 1).Value = "Guadalupe"
$ws.Cells.Item(149, 2).Value = 153
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(149, 4).Value = 104
$ws.Cells.Item(149, 5).Value = 36
$ws.Cells.Item(149, 6).Value = 4
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 13

$ws.Cells.Item(150, 1).Value = "Gibraltar"
$ws.Cells.Item(150, 2).Value = 144
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 141
$ws.Cells.Item(150, 5).Value = 3
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 0

$ws.Cells.Item(151, 1).Value = "Brunei"
$ws.Cells.Item(151, 2).Value = 141
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 132
$ws.Cells.Item(151, 5).Value = 8
$ws.Cells.Item(151, 6).Value = 2
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 1

$ws.Cells.Item(163, 2).Value = 93
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(163, 4).Value = 34
$ws.Cells.Item(163, 5).Value = 49
$ws.Cells.Item(163, 6).Value = 3
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 10

$ws.Cells.Item(205, 1).Value = "Montserrat"
$ws.Cells.Item(205, 2).Value = 11
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 7
$ws.Cells.Item(205, 5).Value = 3
$ws.Cells.Item(205, 6).Value = 1
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 1

$ws.Cells.Item(206, 1).Value = "Seychelles"
$ws.Cells.Item(206, 2).Value = 11
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 8
$ws.Cells.Item(206, 5).Value = 3
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0
